$wb = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item("Questions")
$ws7.Activate()
$excel.ActiveWindow.Zoom = 87
